# Actualización de metricas de kanban
# Insert 7 new "Task" rows into the lead-times table (rows 38-44), fill in
# the previously-empty row 9 (C9:F9), and let the histogram / WIP table
# and charts below shift down accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fill the previously empty row 9 cells (C9:F9) — task "HU9" released on
#    day 6 (row 9 of the Column counts table).
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 7

# 2) Insert 7 rows before row 38 so the new task rows can be added, pushing
#    the "WIP (Finger chart)" section (and everything after it) down by 7.
$ws.Rows("38:44").Insert()

# Carry the bordered table formatting (style used by rows 17-37) down into
# the newly inserted rows before filling them with data.
$ws.Range("B37:F37").Copy()
$ws.Range("B38:F44").PasteSpecial(-4122)

# 3) Populate the 7 new rows (38-44) with the next tasks in the lead-time
#    table, continuing the pattern used by rows 17-37.
$taskRows = @(
    @{ Row = 38; B = 22; C = 9;  D = 10; F = "HU10" },
    @{ Row = 39; B = 23; C = 9;  D = 10; F = "HU6" },
    @{ Row = 40; B = 24; C = 10; D = 10; F = "HU13" },
    @{ Row = 41; B = 25; C = 10; D = 10; F = "HU14" },
    @{ Row = 42; B = 26; C = 10; D = 10; F = "HU19" },
    @{ Row = 43; B = 27; C = 10; D = 10; F = "HU20" },
    @{ Row = 44; B = 28; C = 10; D = 10; F = "HU21" }
)

foreach ($task in $taskRows) {
    $r = $task.Row
    $ws.Range("B$r").Value = $task.B
    $ws.Range("C$r").Value = $task.C
    $ws.Range("D$r").Value = $task.D
    $ws.Range("E$r").Formula = "=D$r-C$r"
    $ws.Range("F$r").Value = $task.F
}

$ws.Range("F45").Select()

$wb.Application.CalculateFullRebuild()
